$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the product_code values for the four existing "replay" rows (5-8).
# The style (text number format) on these cells is preserved automatically
# when only the value is cleared.
$ws.Range("G5").Value = $null
$ws.Range("G6").Value = $null
$ws.Range("G7").Value = $null
$ws.Range("G8").Value = $null

# Add a new row for the "Tokyo Eye Shot" supplement, moved here from another
# (shadowrun city file) sheet.
$ws.Range("A9").Value = 1995
$ws.Range("B9").Value = "TOKYO EYE‐SHOT シャドウランシティファイル"
$ws.Range("C9").Value = "Tokyo Eye Shot: Shadowrun City File"
$ws.Range("D9").Value = "Fujimi Shobo"
$ws.Range("E9").Value = "shadowrun_tokyo_eye_shot.jpg"
$ws.Range("F9").Value = "supplement"
$ws.Range("G9").Value = "12-3"

# Match the style used by the other product_code cells (text number format).
$ws.Range("G9").NumberFormat = "@"

# Update the selected cell to mirror what was saved in the workbook.
$ws.Range("G8").Select()
